$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'بوتولیسم بیماری کدام مورد است؟'
$ws.Range("A3").Value = 'ناحیه ویل شامل همه موارد زیر به جز کدام است؟'
$ws.Range("A4").Value = 'gtr مربوط به کدام یک از موارد زیر است؟'
$ws.Range("A5").Value = 'کدام یک از گزینه‌های زیر مطلوب‌ترین توالی رویش دندان‌های دائمی فک پایین را نشان می‌دهد؟'
$ws.Range("A6").Value = 'در صورت متفاوت بودن اندازه‌های باکس پروگزیمال دندان‌های مجاور، کدام باکس پروگزیمال باید ترجیح داده شود؟'
$ws.Range("A7").Value = 'دستگاه گولژی تاندون پیام را به cns منتقل می‌کند و به چه چیزی وابسته است؟'
$ws.Range("A8").Value = 'آدرنالین به لیدوکائین اضافه می‌شود تا اثر آن را طولانی‌تر کند و جذب آن را در جریان خون کاهش دهد. این نسبت به چه صورت است؟'
$ws.Range("A9").Value = 'حجم خون در یک کودک، به ازای هر کیلوگرم از وزن بدن'
$ws.Range("A10").Value = 'تحت شرایط فیزیولوژیک طبیعی، فشار مایع مغزی-نخاعی (csf) با کدام عامل زیر متناسب است؟'
$ws.Range("A11").Value = 'کودک طبیعی بین چه سنینی دچار ناروانی گفتار می‌شود؟'
$ws.Range("A12").Value = 'یک بازیکن ۲۰ ساله بسکتبال دانشگاهی پس از احساس درد قفسه سینه و تپش قلب در حین تمرین به کلینیک اورژانس دانشگاه آورده می‌شود. تنگی نفس یا تاکیپنه وجود ندارد. وی سابقه خانوادگی بیماری قلبی را رد می‌کند و سابقه اجتماعی وی منفی از نظر مصرف الکل یا مواد مخدر است. معاینه قلبی قابل توجهی ندارد و نوار قلب انقباضات زودرس بطنی (pvcs) مکرر را نشان می‌دهد. کدام یک از گزینه‌های زیر مناسب‌ترین اقدام بعدی در ارزیابی و/یا مدیریت این بیمار است؟'
$ws.Range("A13").Value = 'علامت کوه فوجی ویژگی کدام یک از موارد زیر است؟'
$ws.Range("A14").Value = 'دختر ۱۶ ساله با قد کوتاه، آمنوره و فاصله زیاد بین نوک پستان‌ها. کاریوتایپ چیست؟'
$ws.Range("A15").Value = 'اگر یک فنر به اندازه 0.5 میلی‌متر فعال شود و نیرویی معادل 150 گرم ایجاد کند، نرخ انحراف بار فنر چقدر خواهد بود؟'
$ws.Range("A16").Value = 'شایع ترین تومور مدیاستین خلفی کدام است؟'
$ws.Range("A17").Value = 'سندرم مالوری وایس معمولاً در چه افرادی شایع است؟'
$ws.Range("A18").Value = 'درمان استاندارد توده آپاندیسیتی چیست؟'
$ws.Range("A19").Value = 'میزان طبیعی تولید مایع زلالیه -'
$ws.Range("A20").Value = 'یک بیمار مرد 80 ساله با تومور خط میانی فک پایین مراجعه می‌کند که حاشیه آلوئولار را درگیر کرده است. بیمار بی‌دندان است. درمان انتخابی چیست؟'
$ws.Range("A21").Value = 'مانور هیملیخ برای چه موردی استفاده می‌شود؟'
$ws.Range("A22").Value = 'شایع‌ترین الگوی بالینی کارسینوم سلول بازال چیست؟'
$ws.Range("A23").Value = 'جفت سرراهی با همه موارد زیر مرتبط است به جز:'
$ws.Range("A24").Value = 'چین dennies morgan نشانگر کدام بیماری است؟ مارس 2013'
$ws.Range("A25").Value = 'کدام آزمایش به تفکیک بین استرابیسم همراه و استرابیسم فلجی کمک می‌کند؟'
$ws.Range("A26").Value = 'کدام یک از داروهای خواب‌آور زیر باعث تسهیل اثرات مهاری gaba می‌شود اما فاقد خواص ضد تشنج یا شل‌کنندگی عضلانی است و تأثیر minimal بر معماری خواب دارد؟'
$ws.Range("A27").Value = 'تشخیص افتراقی برای "سردرد رعدآسا" شامل همه موارد زیر است به جز:'
$ws.Range("A28").Value = 'بیمار یافته‌های زیر را دارد- بزرگ‌شدگی طحال، تعداد کم گلبول‌های قرمز، تعداد طبیعی گلبول‌های سفید و پلاکت‌ها و گسترش خون محیطی سلول‌های قطره‌اشکی را نشان می‌دهد. آسپیراسیون مکرر مغز استخوان ناموفق است. تشخیص احتمالی چیست؟'
$ws.Range("A29").Value = 'در نمونه بیوپسی کلیه، رنگ جونز متنامین نقره برای چه چیزی استفاده می‌شود؟'
$ws.Range("A30").Value = 'مکانیسم عمل ضد سرطانی فلورواوراسیل چیست؟'
$ws.Range("A31").Value = 'همه موارد زیر در مورد سطح خلفی-میانی کلیه صحیح هستند، به جز:'
$ws.Range("A32").Value = 'داروی انتخابی برای کره هانتینگتون کدام است؟'
$ws.Range("A33").Value = 'بیشترین مقدار عضله صاف، نسبت به ضخامت دیواره، در کدام قسمت وجود دارد؟'
$ws.Range("A34").Value = 'ژن brca-1 با کدام یک از موارد زیر مرتبط است؟'
$ws.Range("A35").Value = 'آنته‌ورژن رحم توسط چه چیزی حفظ می‌شود؟'
$ws.Range("A36").Value = 'ماده شیمیایی مورد استفاده برای ارزیابی کیفی و کمی الکل در هوای بازدمی چیست؟'
$ws.Range("A37").Value = 'در آلودگی مدفوعی دیده نمی‌شود –'
$ws.Range("A38").Value = 'کدام یک از موارد زیر جزو عضلات دو شکمی نیست؟'
$ws.Range("A39").Value = 'پاسخ برانگیخته بینایی (ver) در تشخیص همه موارد زیر مفید است به جز:'
$ws.Range("A40").Value = 'ظاهر حباب صابون در انتهای تحتانی رادیوس، درمان انتخابی چیست؟'
$ws.Range("A41").Value = 'کدام یک از موارد زیر از ویژگی‌های ترک مواد افیونی است؟'
$ws.Range("A42").Value = 'کدام سم از طریق camp واسطه نمی‌شود؟'
$ws.Range("A43").Value = 'کدام یک از موارد زیر اولین محل ظهور رگور مورتس پس از مرگ است؟'
$ws.Range("A44").Value = 'آميفوستين از همه موارد زير محافظت مي کند به جز:'
$ws.Range("A45").Value = 'ضایعات خال اوتا معمولاً کدام عصب جمجمه‌ای را درگیر می‌کند؟'
$ws.Range("A46").Value = 'کمترین مایع تحریک کننده برای صفاق -'
$ws.Range("A47").Value = 'عوارض جانبی مانند تغییر رنگدانه‌های پوست و ایکتیوز در کدام دارو مشاهده می‌شود؟'
$ws.Range("A48").Value = 'یک زن که توسط همسایه مورد حمله قرار گرفته بود، با شکستگی دندان میانی به بخش اورژانس آورده شد. آسیب به صورت داخلی به دهان گسترش یافته بود و همچنین کوفتگی در هر دو پا وجود داشت. ماهیت آسیب چیست؟'
$ws.Range("A49").Value = 'تغییر شکل قلاب چوپانی یک ویژگی رادیولوژیک مشخصه کدام بیماری است؟'
$ws.Range("A50").Value = 'در افتالموپاتی تیروئیدی، کدام مورد ارتباط ندارد؟'
$ws.Range("A51").Value = 'بیماری بلانت همچنین به عنوان چه چیزی شناخته می‌شود؟'
$ws.Range("A52").Value = 'مرد 35 ساله ای با پانکراتیت حاد مراجعه می کند. مایع ایده آل برای انتخاب کدام است؟'
$ws.Range("A53").Value = 'یک مرد 35 ساله با علائم آشفتگی، گیجی عمومی، سردرگمی، توهم و بیشفعالی خودمختار به اورژانس مراجعه کرد. در تاریخچهگیری، مشخص شد که آخرین مصرف الکل او 5 روز پیش بوده است. کدام یک از گزینههای زیر میتواند در درمان فوری این وضعیت استفاده شود؟'
$ws.Range("A54").Value = 'کودک با ضایعه اریتماتوز غیر قابل بلانچ و برجسته در سمت راست صورت، درمان مناسب چیست؟'
$ws.Range("A55").Value = 'آفازی نامی ناشی از اختلال در کدام ناحیه است؟'
$ws.Range("A56").Value = 'همه موارد زیر در مورد هالوتان صحیح است به جز'
$ws.Range("A57").Value = 'ساییدگی انتهای قدامی دنده‌ها در کدام بیماری مشاهده می‌شود؟'
$ws.Range("A58").Value = 'کدام یک از موارد زیر با فعالیت رانش پیچیده زبان مرتبط نیست؟'
$ws.Range("A59").Value = 'همه موارد زیر واکنش های کربوکسیلاسیون مستقل از بیوتین هستند، به جز'
$ws.Range("A60").Value = 'بهترین نشانگر غرق شدن قبل از مرگ چیست؟'
$ws.Range("A61").Value = 'جمعیت 1000 نفر تحت پوشش کدام گزینه قرار می‌گیرد؟'
$ws.Range("A62").Value = 'بیشترین بازجذب سدیم در کدام قسمت انجام می‌شود؟'
$ws.Range("A63").Value = 'عفونت نایسریا مشخصه کمبود کدام یک از موارد زیر است؟'
$ws.Range("A64").Value = 'اندازه کامپوزیت های میکروفیل در چه محدوده ای قرار دارد؟'
$ws.Range("A65").Value = 'همه موارد زیر در سیستینوری دفع می‌شوند به جز'
$ws.Range("A66").Value = 'مفهوم "ستون ثبات ستون فقرات" توسط چه کسی ارائه شده است؟'
$ws.Range("A67").Value = 'در مرگ‌های ناشی از غواصی با اسکوبا، تمامی بررسی‌های زیر که قبل، حین و پس از کالبدگشایی انجام می‌شوند ارزش بیشتری دارند به جز'
$ws.Range("A68").Value = 'کدام یک از موارد زیر از ویژگی‌های پای چنبری نیست؟'
$ws.Range("A69").Value = 'همه موارد زیر در مورد سیفلیس صحیح است به جز'
$ws.Range("A70").Value = 'ترانسوستیسم چیست؟'
$ws.Range("A71").Value = 'کدام یک از موارد زیر برای مدیریت اولیه پرولاپس رکتوم در کودکان استفاده می‌شود؟'
$ws.Range("A72").Value = 'تمام ژن های زیر در سرطان زایی کولون نقش دارند، به جز:'
$ws.Range("A73").Value = 'بیمار با درد در ران مراجعه می‌کند که با آسپرین تسکین می‌یابد. تصویربرداری با اشعه ایکس یک توده رادیولوسنت احاطه شده توسط اسکلروز را نشان می‌دهد. تشخیص چیست؟'
$ws.Range("A74").Value = 'شایع ترین نوع کارسینومای دهانه رحم: مارس 2007'
$ws.Range("A75").Value = 'کدام یک از موارد زیر در مورد برداشتن استخوان در هنگام کشیدن دندان عقل نهفته صحیح نیست؟'
$ws.Range("A76").Value = 'یک بیمار میانسال به شما ارجاع داده شده است تا گلوکوم زاویه باز را رد کنید. کدام یک از یافته‌های زیر در تشخیص کمک می‌کند؟'
$ws.Range("A77").Value = 'غدد کوپر در کدام قسمت یافت می‌شوند؟'
$ws.Range("A78").Value = 'کدام یک از عوامل زیر ممکن است در بوی بد دهان نقش داشته باشد، زمانی که سطح مخاطی خشک می‌شود؟'
$ws.Range("A79").Value = 'شاخص پین نیتروس اکساید چیست؟'
$ws.Range("A80").Value = 'در مورد کیسه منی کدام گزینه صحیح است؟'
$ws.Range("A81").Value = 'داروی انتخابی برای زن باردار مشکوک به داشتن نوزاد با هیپرپلازی مادرزادی آدرنال کدام است؟'
$ws.Range("A82").Value = 'چه موردی در مورد مافنید استات نادرست است؟'
$ws.Range("A83").Value = 'همه موارد زیر در مورد ایمنی جمعی برای بیماری‌های عفونی صحیح هستند به جز'
$ws.Range("A84").Value = 'تغییر شکل ساعت شنی معده در کدام یک از شرایط زیر مشاهده می‌شود؟'
$ws.Range("A85").Value = 'همه موارد زیر از علل متابولیک بیماری کبدی هستند به جز'
$ws.Range("A86").Value = 'عصب پرونئال عمقی چه نواحی را عصب‌دهی می‌کند؟'
$ws.Range("A87").Value = 'باند کامپوزیت روی دندان زمانی مورد نیاز است که ضخامت باقیمانده عاج دندان چقدر باشد؟'
$ws.Range("A88").Value = 'il-6 به طور مداوم در کدام بیماری افزایش می یابد؟'
$ws.Range("A89").Value = 'شایع ترین نوع کارسینوم برونکوژنیک اولیه کدام است؟'
$ws.Range("A90").Value = 'یک خانم 35 ساله به مدت 4 ماه گذشته قاعدگی نداشته است. سطح سرمی fsh و lh او بالا و استرادیول او پایین است. علت احتمالی این وضعیت چیست؟'
$ws.Range("A91").Value = 'یک ورزشکار 35 ساله با قد 184 سانتی‌متر، فاصله دست‌ها 194 سانتی‌متر، ضربان قلب 64 در دقیقه، فشار خون 148/64 میلی‌متر جیوه در معاینه معمول، سمع قفسه سینه یک سوفل دیاستولیک طولانی در فضای بین دنده‌ای دوم سمت راست نشان می‌دهد. تشخیص احتمالی چیست؟'
$ws.Range("A92").Value = 'اولین ماده بیحس کننده موضعی که به صورت بالینی استفاده شد کدام است؟'
$ws.Range("A93").Value = 'یک زن ۳۵ ساله با ناباروری اولیه مراجعه می‌کند. در معاینه، یک توده در لگن قابل لمس است. سونوگرافی یک ضایعه کیستیک در تخمدان با ظاهر شیشه‌ای زمین‌شده بدون جریان عروقی داخلی نشان می‌دهد. سطح ca-125 او ۹۰ u/ml است. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A94").Value = 'کدام یک از موارد زیر به گروه b آربوویروس تعلق ندارد؟'
$ws.Range("A95").Value = 'بی حس کننده موضعی طبیعی کدام است؟ (1997)'
$ws.Range("A96").Value = 'hla در کجا قرار دارد؟'
$ws.Range("A97").Value = 'کاردیومیوپاتی در کدام یک از موارد زیر وجود ندارد؟'
$ws.Range("A98").Value = 'ویژگی پاتولوژیک در گرانولوماتوز وگنر در بیوپسی کلیه چیست؟'
$ws.Range("A99").Value = 'کدام گزینه در مورد عصب پنجم جمجمه ای نادرست است؟'
$ws.Range("A100").Value = 'شایع ترین وضعیت پیش بدخیمی در کارسینومای دهانی کدام است؟'
$ws.Range("A101").Value = 'یک زن ۲۵ ساله با رنگ پریدگی خفیف و هپاتوسپلنومگالی متوسط مراجعه کرده است. هموگلوبین او ۹۲ گرم در لیتر و سطح هموگلوبین جنینی ۶۵٪ بود. او تاکنون هیچ انتقال خونی دریافت نکرده است. به احتمال زیاد او از چه بیماری رنج می‌برد؟'
$ws.Range("A102").Value = 'یک کودک 2 ساله سابقه چندین شکستگی استخوان با ضربه‌های جزئی دارد. در معاینه، هپاتوسپلنومگالی و فلج اعصاب جمجمه‌ای ii، vii و viii مشاهده می‌شود. آزمایش‌های آزمایشگاهی پان‌سیتوپنی را نشان می‌دهند. رادیوگرافی‌ها استخوان‌های اسکلروتیک منتشر و متقارن با متافیزهای نامنظم را نشان می‌دهند. تجزیه مولکولی استخوان او نقص در تولید آنزیم کربنیک آنهیدراز برای حل کردن کریستال هیدروکسی‌آپاتیت را نشان می‌دهد. او با پیوند سلول‌های بنیادی خونساز درمان می‌شود. کدام یک از سلول‌های زیر در استخوان‌های او به احتمال زیاد از نظر عملکردی کمبود داشته و پس از پیوند جایگزین شده است؟'
$ws.Range("A103").Value = 'شایع‌ترین تومور بدخیم استخوانی کدام است؟'
$ws.Range("A104").Value = 'ناحیه ورنیکه در کجا قرار دارد؟'
$ws.Range("A105").Value = 'تمام موارد زیر در مورد زخم تیفوئید صحیح است به جز کدام یک؟'
$ws.Range("A106").Value = 'بیمار با درد در چشم چپ همراه با اختلال بینایی مراجعه کرده است، همچنین سابقه ضربه غیرنافذ به چشم را 4 ماه قبل دارد. اولین بررسی انتخابی چیست؟'
$ws.Range("A107").Value = 'رنگ مورد استفاده در میکروسکوپی فلورسنت -'
$ws.Range("A108").Value = 'نوزاد با وزن بسیار کم هنگام تولد به چه نوزادی گفته می‌شود؟'
$ws.Range("A109").Value = 'شایع‌ترین موقعیت آپاندیس کدام است؟'
$ws.Range("A110").Value = 'مرگ ناگهانی پس از شستشوی سینوس فک بالا به چه علتی رخ می‌دهد؟ dnb 10'
$ws.Range("A111").Value = '"herald patch" در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A112").Value = 'میتومایسین-سی موضعی در چه مواردی استفاده می‌شود؟'
$ws.Range("A113").Value = 'bcg چیست؟'
$ws.Range("A114").Value = 'همه موارد زیر در مورد siadh صحیح هستند، به جز -'
$ws.Range("A115").Value = 'کدام گزینه در مورد کارسینوئید آپاندیس صحیح است؟'
$ws.Range("A116").Value = 'ظاهر هایپر دانس دو محدب در سی تی اسکن جمجمه در کدام مورد دیده می‌شود؟'
$ws.Range("A117").Value = 'سندرم آدرنوژنیتال بیشتر به دلیل چه چیزی ایجاد می‌شود؟'
$ws.Range("A118").Value = 'کدام سلول‌ها نسبت به هیپوکسی حساس‌تر هستند؟'
$ws.Range("A119").Value = 'برآمدگی محوری چشم (آکسیال پروپتوزیس) توسط تومورهای واقع در کدام ناحیه ایجاد می‌شود؟'
$ws.Range("A120").Value = 'کدام گزینه در مورد رفلکس عقب‌کشیدن (فلکشن) صحیح نیست؟'
$ws.Range("A121").Value = 'عبارت صحیح در مورد بوپیواکائین کدام است؟'
$ws.Range("A122").Value = '"la facies sympathique" شرایطی است که در موارد زیر دیده می‌شود:'
$ws.Range("A123").Value = 'کدام یک از موارد زیر جزو جراحی‌های مورد استفاده در مدیریت بیماری منیر نیست؟'
$ws.Range("A124").Value = 'فشار خون در بطن راست'
$ws.Range("A125").Value = 'پدیده سوختگی گردنی ناشی از چیست؟'
$ws.Range("A126").Value = 'یک مرد ۶۰ ساله پس از سقوط ناگهانی در توالت، فشار خون ۹۰/۵۰ میلی‌متر جیوه و نبض ۱۰۰ در دقیقه داشت. بستگان او گزارش دادند که مدفوعش سیاه/تیره رنگ بوده است. تاریخچه دقیق‌تر نشان داد که او سابقه فشار خون بالا و بیماری عروق کرونر دارد و به‌طور منظم آسپرین، آتنولول و سوربیترات مصرف می‌کرده است. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A127").Value = 'کدام یک از موارد زیر به طور ترجیحی پلاسمینوژن متصل به فیبرین را فعال می‌کند و از حالت لیتیک سیستمیک جلوگیری می‌کند:'
$ws.Range("A128").Value = 'تمامی این ساختارها در دیواره جانبی بینی یافت می‌شوند به جز -'
$ws.Range("A129").Value = 'بیماری سل در عفونت hiv باعث ایجاد بیماری ریوی شبیه به بیماری پس از مرحله اولیه در افراد عادی می‌شود زمانی که در ........................ مرحله عفونت hiv رخ می‌دهد -'
$ws.Range("A130").Value = 'در مورد کارسینوم فیبرولاملار، کدام گزینه صحیح است؟'
$ws.Range("A131").Value = 'داروی انتخابی برای مالاریای مقاوم به کلروکین ناشی از پلاسمودیوم فالسیپاروم در گروه سنی کودکان کدام است؟'
$ws.Range("A132").Value = 'بینی زین اسبی چیست؟'
$ws.Range("A133").Value = 'یک زن ۲۶ ساله در روز چهارم پس از زایمان یک نوزاد پسر ۴.۵ کیلوگرمی مراجعه می‌کند. او نگران پای راست خود است که از زمان زایمان بی‌حس و ضعیف شده است. راه رفتن برای او دشوار شده زیرا پای راستش تمایل به افتادن دارد و انگشتان پایش کشیده می‌شوند. وقتی درباره روند زایمان از او سؤال می‌شود، گزارش می‌دهد که اپیدورال با تسکین درد رضایت‌بخش داشته اما مرحله فشار زایمان دشوار و طولانی (۳ ساعت) بوده است. او در طول زایمان روی تخت قرار گرفته بود. او هرگونه کمردرد یا مشکل در پای دیگر را رد می‌کند. در معاینه، کاهش حس در روی پای راست و طرف خارجی ساق پا همراه با ناتوانی در دورسی فلکشن پای راست که منجر به افتادگی پا شده است، مشاهده می‌شود. ادم محیطی minimal در هر دو اندام تحتانی دیده می‌شود. کدام عصب به احتمال زیاد تحت فشار قرار گرفته است؟'
$ws.Range("A134").Value = 'دیورز در کدام مرحله از آسیب حاد کلیوی مشاهده می‌شود؟'
$ws.Range("A135").Value = 'رشد استالاسایت توسط کدام عامل است؟'
$ws.Range("A136").Value = 'شایع‌ترین علت واکنش‌های تزریق خون کشنده چیست؟'
$ws.Range("A137").Value = 'کدام یک از عوامل بیهوشی زیر باعث تحریک هیپرترمی بدخیم نمی‌شود؟'
$ws.Range("A138").Value = 'ماده بخیه روده توسط چه فرآیندی جذب می‌شود؟'
$ws.Range("A139").Value = 'آنکوژن مرتبط با لنفوم بورکیت کدام است؟'
$ws.Range("A140").Value = 'در سوختگی‌های شدید، کمترین مورد استفاده کدام است؟'
$ws.Range("A141").Value = 'کدام گزینه در مورد استئوژنز ایمپرفکتا صحیح نیست؟'
$ws.Range("A142").Value = 'ترکیب ورنی چیست؟'
$ws.Range("A143").Value = 'یک زن ۶۲ ساله با توده پستانی که ۶ روز پیش کشف کرده است مراجعه می‌کند. بیوپسی پستان کارسینوم لوبولار درجا را نشان می‌دهد. در مقایسه با سلول‌های اپیتلیال طبیعی لوبول پستان، این سلول‌های بدخیم به احتمال زیاد کاهش بیان کدام یک از پروتئین‌های زیر را نشان می‌دهند؟'
$ws.Range("A144").Value = 'عصب دهی پاراسمپاتیک به غدد بزاقی توسط کدام اعصاب انجام می‌شود؟'
$ws.Range("A145").Value = 'کدام یک از داروهای زیر در درمان بیماری پارکینسون، آنتاگونیست nmda است؟'
$ws.Range("A146").Value = 'مدیریت مناسب برای بیمار با زخم تمیز روی ساعد که کمی از دست دادن بافت دارد چیست؟ او ۱۲ سال پیش واکسن کزاز (تتانوس توکسوئید) دریافت کرده است.'
$ws.Range("A147").Value = 'در آزمایش ادرار بیمار مبتلا به هماچوری و گلومرولونفروز چه چیزی مشاهده می‌شود؟'
$ws.Range("A148").Value = 'تومور ثانویه می‌تواند از طریق تمام موارد زیر به مدار چشم گسترش یابد به جز'
$ws.Range("A149").Value = 'همه موارد زیر در مورد اشرشیاکلی انتروهموراژیک صحیح است به جز؟'
$ws.Range("A150").Value = 'دیپلوپی (دو بینی) بیشتر در کدام مورد شایع است؟'
$ws.Range("A151").Value = 'شایع‌ترین علت بزرگ‌شدن سایه قلب در عکس‌برداری اشعه ایکس از یک کودک چیست؟'
$ws.Range("A152").Value = 'اگرچه گلیکولیز و مسیر فسفات پنتوز چندین متابولیت مشترک دارند، اما در موارد دیگر به طور قابل توجهی متفاوت هستند. تمام عبارات زیر صحیح هستند، به جز:'
$ws.Range("A153").Value = 'لیگامان‌های مفصل شانه همه موارد زیر هستند به جز'
$ws.Range("A154").Value = 'یک مرد 49 ساله با علائمی که پس از یک آخر هفته طولانی مصرف زیاد الکل ایجاد شده است، مراجعه می‌کند. سطح سرمی گاما-گلوتامیل ترانسفراز (ggt) او 65 iu/l است. بیوپسی کبد تغییرات چربی (استئاتوز) در تعداد زیادی از هپاتوسیت‌ها را نشان می‌دهد. این ناهنجاری کبدی در این بیمار به احتمال زیاد ناشی از کدام مورد است؟'
$ws.Range("A155").Value = 'یک مرد 85 ساله تحت عمل برداشت ترانس اورترال پروستات قرار گرفت. بررسی هیستولوژیک نمونه او کانون‌های آدنوکارسینوم را نشان داد. مدیریت این مورد چه خواهد بود؟'
$ws.Range("A156").Value = 'رژیم‌های غذایی با پتاسیم بالا شامل کدام‌یک از موارد زیر است؟ الف) سبزیجات برگ‌دار سبز ب) شیر ج) موز د) نمک یددار ه) ماهی'
$ws.Range("A157").Value = 'تمامی گزینه‌های زیر صحیح هستند به جز:'
$ws.Range("A158").Value = 'پس از ۳ هفته از تشکیل کیست کاذب پانکراس به اندازه ۵ سانتی‌متر، کدام روش برای مدیریت آن مناسب است؟'
$ws.Range("A159").Value = 'اگر به dna نمک اضافه شود چه اتفاقی می‌افتد؟ (pgi dec 2008)'
$ws.Range("A160").Value = 'ویژگی‌های یک کاندیدای ایده‌آل برای قرار دادن آی‌یو‌دی مسی شامل همه موارد زیر به جز کدام است؟'
$ws.Range("A161").Value = 'آنوری به عنوان دفع ادرار کمتر از چه مقدار تعریف می‌شود؟'
$ws.Range("A162").Value = 'یک نوزاد 12 ساعت پس از تولد مکونیوم سیاه رنگ دفع می‌کند. عبارت صحیح کدام است؟'
$ws.Range("A163").Value = 'لیزین در کدام یک از موارد زیر وجود ندارد؟'
$ws.Range("A164").Value = 'برجستگی مجرای ادرار (urethral crest) در اثر کدام عامل در مجرای ادرار دیده می‌شود؟'
$ws.Range("A165").Value = 'محتوای ors سازمان جهانی بهداشت شامل کدام یک است؟'
$ws.Range("A166").Value = 'در پزشکی قانونی، شاخص سفالیک (cephalic index) در تعیین کدام مورد استفاده می‌شود؟'
$ws.Range("A167").Value = 'اختلال در "وضعیت آنتیسیپاتوری" در کدام یک از موارد زیر مشاهده میشود؟'
$ws.Range("A168").Value = '"هدایت به رویش" اصطلاح دیگری است برای'
$ws.Range("A169").Value = 'تشخیص پاتولوژی چشمی را انجام دهید.'
$ws.Range("A170").Value = 'کدام یک از ابزارهای زیر برای کورتاژ ناحیه فورکیشن استفاده می‌شود؟'
$ws.Range("A171").Value = 'یک سیکلوپلژی قوی کدام است؟'
$ws.Range("A172").Value = 'همه موارد زیر درباره n20 نادرست هستند به جز؟'
$ws.Range("A173").Value = 'در مورد سندرم درسلر همه موارد زیر صحیح است، به جز'
$ws.Range("A174").Value = 'کدام یک از عبارات زیر در مورد اجسام هاسال صحیح است؟'
$ws.Range("A175").Value = 'یک کودک ۴ ساله با یک توده بدون درد منفرد در گردن مراجعه کرده است. والدین گزارش می‌دهند که توده معمولاً در زمان عفونت دستگاه تنفسی فوقانی متورم و حساس می‌شود. در معاینه، یک توده صاف، بدون حساسیت و نوسانی در امتداد یک سوم پایینی مرز قدامی-میانی عضله استرنوکلیدوماستوئید بین عضله و پوست رویی مشاهده شد. سی‌تی اسکن کنتراست گردن انجام شد. این ساختار به دلیل عدم انسداد کدام یک از موارد زیر ایجاد می‌شود؟'
$ws.Range("A176").Value = 'کویینود سگمان نهم کبد متعلق به کدام بخش است؟'
$ws.Range("A177").Value = 'درماتوپاتی ناشی از کم کاری تیروئید به صورت زیر ظاهر می شود:'
$ws.Range("A178").Value = 'کدام یک از موارد زیر در کولیت اولسراتیو (u.c.) بیشترین درگیری را دارد؟'
$ws.Range("A179").Value = 'کدام یک از موارد زیر یک مخمر واقعی است؟'
$ws.Range("A180").Value = 'رئوباز نشان دهنده چیست؟'
$ws.Range("A181").Value = 'یک خانم جوان با سابقه اختلال میلوپرولیفراتیو، با شکایت از ناراحتی پیشرونده شکمی و آسیت به مدت ۳ روز مراجعه کرده است. کدام یک از تشخیص‌های زیر محتمل‌تر است؟'
$ws.Range("A182").Value = 'شایع ترین علت اوتیت مدیا حاد در کودکان چیست؟'
$ws.Range("A183").Value = 'شایع‌ترین عفونت دهانی در دیابت شیرین چیست؟'
$ws.Range("A184").Value = 'کدام بخش از گوش منشأ آن از هر سه لایه زاینده است؟'
$ws.Range("A185").Value = 'تمامی عبارات مربوط به نرخ خام تولد صحیح هستند، به جز:'
$ws.Range("A186").Value = 'آزمون چالش اکسی توسین برای ارزیابی سلامت جنین در همه موارد زیر منع مصرف دارد به جز:'
$ws.Range("A187").Value = 'لوکوریا در تمام موارد زیر دیده می‌شود به جز-'
$ws.Range("A188").Value = 'در مورد اجسام پساموما همه موارد زیر صحیح است به جز'
$ws.Range("A189").Value = 'همه موارد زیر از c-amp به عنوان پیام‌رسان دوم استفاده می‌کنند به جز:'
$ws.Range("A190").Value = 'تنیس البو با کدام مورد مشخص می‌شود؟'
$ws.Range("A191").Value = 'التهاب بهاری چشم ممکن است با کدام یک از موارد زیر همراه باشد؟'
$ws.Range("A192").Value = 'همه موارد زیر می‌توانند باعث خفگی تروماتیک شوند، به جز:'
$ws.Range("A193").Value = 'کدام یک از عوارض زیر به احتمال زیاد پس از دریافت چند واحد انتقال خون رخ می‌دهد؟'
$ws.Range("A194").Value = 'شکل اصلی ذخیره و گردش ویتامین d کدام است؟'
$ws.Range("A195").Value = 'یک زن ۲۷ ساله با سابقه خانوادگی بیماری خودایمنی، با شکایت از بثورات پوستی و دردهای مفصلی مکرر ۳ ماه پس از زایمان مراجعه کرده است. به احتمال زیاد کدام یک از اختلالات زیر را دارد؟'
$ws.Range("A196").Value = 'تریاد شارکو در کدام مورد مشاهده می‌شود؟ سپتامبر 2006'
$ws.Range("A197").Value = 'دروازه‌بان قلب چیست؟'
$ws.Range("A198").Value = 'نوع سی‌تی اسکن مورد استفاده برای تعیین ترکیب شیمیایی سنگ‌های کلیه چیست؟'
$ws.Range("A199").Value = 'نمره آپگار شامل همه موارد زیر است به جز -'
$ws.Range("A200").Value = 'مجرای اندولنفاتیک به کجا تخلیه می‌شود؟'
$ws.Range("A201").Value = 'از کل مقدار بزاق ترشح شده توسط تمام غدد بزاقی در یک روز، بیش از نیمی از آن توسط کدام غده ترشح می‌شود؟'
$ws.Range("A202").Value = 'دوز دیگوکسین در کدام مورد تغییر نمی‌کند؟'
$ws.Range("A203").Value = 'کدام یک از تومورهای زیر به پرتودرمانی حساس است؟'
$ws.Range("A204").Value = 'کوتاه‌اثرترین شل‌کننده عضلات اسکلتی کدام است؟'
$ws.Range("A205").Value = 'تعریف هیپوگلیسمی در نوزادان چیست؟'
$ws.Range("A206").Value = 'یک مرد ۲۵ ساله پس از تصادف رانندگی با فشار خون ۱۰۰/۸۰ میلی‌متر جیوه و نبض ۸۴ در دقیقه مراجعه می‌کند. بهترین مایع برای احیا چیست؟'
$ws.Range("A207").Value = 'یک زن ۴۵ ساله با درد و سفتی گردن به مدت ۳ ماه مراجعه کرده است. تشخیص به صورت زیر است. کدام یک از موارد زیر جزو درمان‌های جراحی این بیماری نیست؟'
$ws.Range("A208").Value = 'ضدعفونی کردن آب لازم نیست اگر منبع آن از کدام مورد باشد؟'
$ws.Range("A209").Value = 'یک پسر ۳ ساله از دوران نوزادی تأخیر رشدی پیشرونده، آتاکسی، تشنج و خنده نامناسب داشته است. کودک کاریوتایپ طبیعی ۴۶، xy دارد، اما تجزیه و تحلیل dna نشان می‌دهد که هر دو کروموزوم شماره ۱۵ خود را از پدر به ارث برده است. این یافته‌ها به احتمال زیاد نشانگر کدام یک از مکانیسم‌های ژنتیکی زیر است؟'
$ws.Range("A210").Value = 'کودکی با یک لکه سفید روی لوزه‌ها مراجعه می‌کند، بهترین روش تشخیصی کشت در'
$ws.Range("A211").Value = 'کدام یک از داروهای زیر با روش استنشاقی در درمان آنژین استفاده شده و دارای شروع اثر بسیار سریع و مدت اثر کوتاه است؟'
$ws.Range("A212").Value = 'زاویه بزرگ کاپا باعث کدام مورد می‌شود؟'
$ws.Range("A213").Value = 'ویژگی‌های رادیولوژیک در راشیتیسم کدامند؟'
$ws.Range("A214").Value = 'پس از انزال در واژن، اسپرم تا چه مدت قابلیت حرکت خود را حفظ می‌کند؟'
$ws.Range("A215").Value = 'بیماری فوربس به دلیل کمبود کدام آنزیم ایجاد می‌شود؟'
$ws.Range("A216").Value = 'عبارت نادرست در مورد تومور کام سخت کدام است؟'
$ws.Range("A217").Value = 'هنگام انجام dcr، استئوم در سطح کدام شاخک ایجاد می‌شود؟'
$ws.Range("A218").Value = 'یک زن 80 ساله پس از زمین خوردن دچار عدم توانایی در راه رفتن شده و در معاینه، دفورمیتی چرخش خارجی وجود دارد. slr امکان‌پذیر نیست و پهن شدن تروکانتر مشاهده می‌شود. تشخیص احتمالی چیست؟'
$ws.Range("A219").Value = 'شایع‌ترین محل پارگی در سندرم بوئرهاو کدام است؟'
$ws.Range("A220").Value = 'آگونیست جزئی می‌تواند اثرات آگونیست کامل را آنتاگونیست کند زیرا دارای'
$ws.Range("A221").Value = 'یک کودک 10 ماهه با صدای غیرطبیعی در دم مراجعه می‌کند؛ مادر شکایت دارد که این صدا هنگام گریه کودک افزایش می‌یابد و زمانی که کودک در حالت دمر می‌خوابد کاهش می‌یابد. تشخیص احتمالی چیست؟'
$ws.Range("A222").Value = 'طول عمر گلبول های قرمز نوزادان چقدر است؟'
$ws.Range("A223").Value = 'کدام یک از مواد ترمیمی زیر برای پوسیدگی سطح ریشه توصیه می‌شود؟'
$ws.Range("A224").Value = 'کدام یک از موارد زیر از ویژگی‌های کم‌خونی مگالوبلاستیک نیست؟ سپتامبر 2010'
$ws.Range("A225").Value = 'مس جزئی از کدام آنزیم است؟'
$ws.Range("A226").Value = 'ضایعات متعدد پسوریازیس روی دست ها. درمان انتخابی چیست؟'
$ws.Range("A227").Value = 'در میلوم متعدد موارد زیر مشاهده می‌شود -('
$ws.Range("A228").Value = 'کدام یک از بدخیمی‌های انسانی زیر بسیار بدخیم، با پیشرفت سریع است و تأخیر حتی 1-2 روزه در درمان می‌تواند منجر به مرگ شود؟'
$ws.Range("A229").Value = 'کدام یک از ساختارهای زیر از تلانسفالون مشتق شده است؟'
$ws.Range("A230").Value = 'افزایش ایزوله aptt در کدام مورد مشاهده می‌شود؟'
$ws.Range("A231").Value = 'محتمل‌ترین بدخیمی که در موارد گواتر طولانی‌مدت ایجاد می‌شود کدام است؟'
$ws.Range("A232").Value = 'محتوای سدیم (na+) در رینگر لاکتات ................. میلی‌اکیوالان در لیتر است.'
$ws.Range("A233").Value = 'کدام یک از رویدادهای بلوغ در دختران وابسته به استروژن نیست؟'
$ws.Range("A234").Value = 'همه موارد زیر به جز یکی از عوارض سرخجه مادرزادی هستند:'
$ws.Range("A235").Value = 'بوی طناب سوخته ناشی از مسمومیت با: aiims 14'
$ws.Range("A236").Value = 'گودال‌های هربرت در کدام بیماری دیده می‌شوند؟'
$ws.Range("A237").Value = 'طغیان تنفسی مربوط به کدام مورد است؟'
$ws.Range("A238").Value = 'کدام یک از موارد زیر یک آنتی سایکوتیک آتیپیکال (غیرمعمول) محسوب نمی‌شود؟'
$ws.Range("A239").Value = 'کدام یک از موارد زیر در درمان سیستیت در دوران بارداری استفاده نمی‌شود؟'
$ws.Range("A240").Value = 'همه موارد زیر در متابولیسم ویتامین d نقش دارند به جز'
$ws.Range("A241").Value = 'فیلم‌های اسکرین با فیلم‌های غیراسکرین چه تفاوتی دارند؟'
$ws.Range("A242").Value = 'سه‌گانه‌ی دفورمیتی مفصل زانو در کدام مورد مشاهده می‌شود؟'
$ws.Range("A243").Value = 'تیامین پیروفسفات یک کوآنزیم است که برای همه موارد زیر مورد نیاز است به جز'
$ws.Range("A244").Value = 'یک مرد 35 ساله با سابقه ترشح مجرای ادرار در سه روز گذشته مراجعه می‌کند. رنگ‌آمیزی گرم ترشحات در زیر نشان داده شده است. کدام یک از موارد زیر در مورد عامل احتمالی صحیح است؟'
$ws.Range("A245").Value = 'ترشح خونی از گوش بیشتر به چه علتی است؟'
$ws.Range("A246").Value = 'هنگامی که یک فرد سالم دچار پنوموتوراکس خودبه‌خودی در ریه راست می‌شود، انتظار می‌رود کدام یک از موارد زیر رخ دهد؟'
$ws.Range("A247").Value = 'سندرم فِلتی با کدام بیماری مرتبط است؟'
$ws.Range("A248").Value = 'همه موارد زیر مشتقات مجرای پارامزونفریک هستند به جز:'
$ws.Range("A249").Value = 'واکوئله شدن زیر هسته ای نشانه مشخص کدام مرحله است؟'
$ws.Range("A250").Value = 'آزمایش غیرتهاجمی برای هلیکوباکتر پیلوری - الف) تست سریع اوره آز ب) تست تنفسی اوره آز ج) آزمایش آنتی ژن مدفوع د) کشت آسپیراسیون معده ه) بیوپسی'
$ws.Range("A251").Value = 'ساختاری که منطقه آب‌‌خیز (watershed) را تشکیل نمی‌دهد:'
$ws.Range("A252").Value = 'تب استخوان شکن'
$ws.Range("A253").Value = 'استخوان مرده در تصویر اشعه ایکس چگونه است؟'
$ws.Range("A254").Value = 'کدام یک از موارد زیر در مورد صدای چهارم قلب (s4) صحیح است؟'
$ws.Range("A255").Value = 'در نوریت جذامی، کدام یک از موارد زیر صحیح است؟  
الف) ضخیم شدن داخل پوستی عصب  
ب) فلج صورت  
ج) در مناطق گرم و مرطوب رخ می‌دهد  
د) اعصاب محیطی قابل لمس هستند  
ه) ممکن است بیش از ۱۰ ضایعه یافت شود'
$ws.Range("A256").Value = 'کدام یک از موارد زیر نمونه‌ای از ردیاب درون دهانی است؟'
$ws.Range("A257").Value = 'بیماران دیالیزی مستعد ابتلا به کدام مسمومیت هستند؟'
$ws.Range("A258").Value = 'ترکیب تقریبی پروتئین‌ها در غشای سلولی چقدر است؟'
$ws.Range("A259").Value = 'کدام یک از موارد زیر در یافته های مایع مغزی نخاعی در مننژیت سلی دیده نمی شود؟'
$ws.Range("A260").Value = 'سورگوم حاوی مقدار خوبی از کدام ماده است؟'
$ws.Range("A261").Value = 'در درمان مزمن فشار خون بالا، افت فشار خون وضعیتی بیشترین میزان را با کدام دارو دارد؟'
$ws.Range("A262").Value = 'یک مرد 55 ساله از ضعف شدید عضلانی و افتادگی پلک ها شکایت دارد. او بیان می کند که علائمش با حرکات تکراری بدتر می شود اما پس از استراحت کوتاه بهبود می یابد. عکس برداری قفسه سینه یک توده مدیاستین قدامی را نشان می دهد. بیوپسی این توده به احتمال زیاد کدام یک از تغییرات پاتولوژیک زیر را نشان می دهد؟'
$ws.Range("A263").Value = 'کدام یک از درگیری‌های اندام زیر معمولاً همراه با درگیری سیستم عصبی مرکزی (cns) در یک بیمار مبتلا به لنفوم دیفوز سلول b بزرگ (dlbcl) مشاهده نمی‌شود؟'
$ws.Range("A264").Value = 'همه موارد زیر در مورد انفارکتوس سفید صحیح است به جز؟'
$ws.Range("A265").Value = 'ویژگی(های) نفروز لیپوئید کدام است؟  
الف) طبیعی در میکروسکوپ نوری  
ب) رسوب اپیتلیال  
ج) اسکلروز گلومرولار  
د) از بین رفتن یکنواخت و منتشر پایه‌های پودوسیت  
ه) اسکلروز توبولار'
$ws.Range("A266").Value = 'اوتیت خارجی بدخیم ناشی از عفونت کدام یک از ارگانیسم های زیر است؟'
$ws.Range("A267").Value = 'تعداد طبیعی پلاکت در کدام یک از موارد زیر دیده می‌شود؟ (pgi nov 2010)'
$ws.Range("A268").Value = 'مطالعه مورد-شاهدی چیست؟'
$ws.Range("A269").Value = 'هوموسیستئین با کدام یک از موارد زیر ارتباط ندارد؟'
$ws.Range("A270").Value = 'کدام یک از کلستریدیاهای زیر غیرتهاجمی است؟'
$ws.Range("A271").Value = 'تنسورهای تارهای صوتی کدامند؟'
$ws.Range("A272").Value = 'کدام یک از گلوکوکورتیکوئیدهای زیر را نمی توان به صورت استنشاقی تجویز کرد؟'
$ws.Range("A273").Value = 'مهم‌ترین سلول‌ها در حساسیت نوع اول کدامند؟'
$ws.Range("A274").Value = 'یک بیمار زن به نام ناندینی با عفونت دستگاه تنفسی فوقانی مراجعه می‌کند. دو روز بعد، او دچار هماچوری می‌شود. تشخیص احتمالی چیست؟'
$ws.Range("A275").Value = 'گیرنده کدام یک از موارد زیر به صورت درون سلولی وجود دارد؟'
$ws.Range("A276").Value = 'ریشه عصبی عصب زیر بغل (اکسیلاری) کدام است؟'
$ws.Range("A277").Value = 'غشای پایه اطراف سلول‌های شوان حاوی کدام یک از کلاژن‌های زیر است؟'
$ws.Range("A278").Value = 'پارکینسونیسم با پارکینسونیسم غیرمعمول در چه موردی تفاوت دارد؟'
$ws.Range("A279").Value = 'یک مرد 48 ساله تحت عمل جراحی برای زخم مزمن دوازدهه قرار می‌گیرد. این روش شامل واگوتومی تنه‌ای و کدام یک از موارد زیر است؟'
$ws.Range("A280").Value = 'کوتاهی قد به عنوان قدی تعریف می‌شود که زیر چه مقداری باشد؟'
$ws.Range("A281").Value = 'استافیلومای خلفی بیشتر در کدام مورد مشاهده می‌شود؟ مارس 2012'
$ws.Range("A282").Value = 'کدام یک از موارد زیر جزو سایه مدیاستینال در سمت راست نیست؟'
$ws.Range("A283").Value = 'یک بیمار مرد در حین انجام معاینات و بررسی‌های معمول سلامت، به صورت اتفاقی در معاینه بالینی یک توده تشخیص داده شد. پزشک سونوگرافی را توصیه کرد و پس از آن بیوپسی با هدایت سونوگرافی انجام شد. بررسی هیستوپاتولوژیک (hpe) تشخیص زیر را نشان داد. بیمار تحت شیمی‌درمانی قرار گرفت که شامل دارویی به دست آمده از گیاه داده شده بود. محل اثر داروی ذکر شده در بالا کجاست؟'
$ws.Range("A284").Value = 'شایع‌ترین عارضه تنگی زیر دریچه‌ای آئورت کدام است؟'
$ws.Range("A285").Value = 'ماتریکس موقت از چه چیزی تشکیل شده است؟'
$ws.Range("A286").Value = 'سرطان پروستات معمولاً به مهره‌ها متاستاز می‌دهد:'
$ws.Range("A287").Value = 'رسوب کمپلکس ایمنی، آسیب‌شناسی پایه‌ای در کدام نوع از واکنش‌های حساسیت است؟'
$ws.Range("A288").Value = 'رایج‌ترین روش جراحی برای زخم دوازدهه سوراخ‌شده کدام است؟'
$ws.Range("A289").Value = 'کدام یک از موارد زیر برای کشش اسکلتی استفاده نمی‌شود؟'
$ws.Range("A290").Value = 'اجسام نگری مشخصه کدام بیماری هستند؟ سپتامبر 2008، مارس 2013'
$ws.Range("A291").Value = 'داروی ضد صرع انتخابی در نوروپاتی محیطی پس از دیابت کدام است؟'
$ws.Range("A292").Value = 'شایع ترین واکنش بیوتغییر فاز اول کدام است؟'
$ws.Range("A293").Value = 'ca-125 در همه موارد زیر افزایش می‌یابد به جز:'
$ws.Range("A294").Value = 'عمل پلیکیشن رترکتورهای تحتانی پلک در کدام مورد نشان داده شده است؟'
$ws.Range("A295").Value = 'اکروزنیس به دلیل کمبود کدام آنزیم ایجاد می‌شود؟'
$ws.Range("A296").Value = 'سوزاندن (زباله‌سوزی) برای کدام یک از موارد زیر انجام می‌شود؟ سپتامبر 2008'
$ws.Range("A297").Value = 'در سنکوپ نوروکاردیوژنیک، کم‌ترین بررسی مفید کدام است؟'
$ws.Range("A298").Value = 'رنگدانه مالاریا توسط چه چیزی تشکیل می‌شود؟'
$ws.Range("A299").Value = 'علامت پای عنکبوتی در ivp نشانگر کدام مورد است؟ مارس 2010'
$ws.Range("A300").Value = 'عوامل خطر برای تغییر بدخیم در یک بیمار بدون علامت با پولیپ کیسه صفرا در سونوگرافی شامل همه موارد زیر به جز کدام است؟'
$ws.Range("A301").Value = 'یک کودک 13 ساله با شکایت از عدم شروع قاعدگی و کاریوتایپ 46xx به درمانگاه زنان مراجعه می‌کند. در معاینه، کلیترومگالی مشاهده می‌شود. کمبود کدام آنزیم در این شرایط محتمل است؟'
$ws.Range("A302").Value = 'هذیان‌ها به احتمال زیاد در کدام مورد دیده نمی‌شوند؟'
$ws.Range("A303").Value = 'یک مرد ۲۰ ساله دارای توده بیضه است که پس از ارکیکتومی، یک تومور بدخیم با تمایز کیسه زرده نشان می‌دهد. کدام یک از نشانگرهای تومور زیر احتمالاً برای پایش بیمار از نظر بیماری عودکننده یا متاستاتیک مفیدتر است؟'
$ws.Range("A304").Value = 'داروی مورد استفاده در تشنج های غیابی که طیف باریکی از فعالیت ضد صرع دارد کدام است؟'
$ws.Range("A305").Value = 'کدام یک از موارد زیر در مورد ناقل تیفوئید صحیح نیست؟'
$ws.Range("A306").Value = 'یک پسر ۱۵ ساله با سابقه یک روزه خونریزی لثه، خونریزی زیر ملتحمه و راش پورپوریک مراجعه کرده است. بررسی‌ها نتایج زیر را نشان داد: هموگلوبین ۶.۴ گرم در دسی‌لیتر، گلبول‌های سفید ۲۶۵۰۰ در میلی‌متر مکعب، پلاکت ۳۵۰۰۰ در میلی‌متر مکعب، زمان پروترومبین ۲۰ ثانیه با کنترل ۱۳ ثانیه، زمان ترومبوپلاستین جزئی ۵۰ ثانیه و فیبرینوژن ۱۰ میلی‌گرم در دسی‌لیتر. اسمیر محیطی نشان‌دهنده لوسمی میلوبلاستیک حاد بود. کدام یک از گزینه‌های زیر محتمل‌ترین تشخیص است؟'
$ws.Range("A307").Value = 'نکروز آواسکولار (avn) سر استخوان ران در کدام مورد شایع‌تر است؟'
$ws.Range("A308").Value = 'ترتیب ظهور ویژگی‌های جنسی ثانویه در مردان –'
$ws.Range("A309").Value = 'متاستازهای گردنی بدون وجود تومور اولیه واضح، اغلب ناشی از کدام مورد است؟'
$ws.Range("A310").Value = 'بیماری ماری و سینتون همچنین به عنوان چه چیزی شناخته می‌شود؟'
$ws.Range("A311").Value = 'مقیاس آپگار در نوزادان چه مواردی را اندازه‌گیری می‌کند؟'
$ws.Range("A312").Value = 'کدام آنتی بادی در مادر مبتلا به لوپوس اریتماتوز سیستمیک (sle) مسئول بیماری مادرزادی قلب در کودک است؟'
$ws.Range("A313").Value = 'سابقه مصرف شیر برنج، کودک با شکم برآمده، آلبومین پایین اما بدون پروتئینوری، تشخیص احتمالی چیست؟'
$ws.Range("A314").Value = 'کدام یک از آنزیم های زیر به عنوان داروی ضد سرطان استفاده می شود؟'
$ws.Range("A315").Value = 'کدام یک از گزینه‌های زیر صحیح است:'
$ws.Range("A316").Value = 'در محاسبه سال‌های عمر از دست رفته بالقوه (ypll)، مخرج کسر چیست؟'
$ws.Range("A317").Value = 'در کمبود نیاسین، همه موارد زیر مشاهده می‌شوند به جز -'
$ws.Range("A318").Value = 'خونریزی زیر ملتحمه در موارد زیر مشهود است به جز؟'
$ws.Range("A319").Value = 'کدام باسیلوس برای آزمایش اثربخشی استریلیزاسیون توسط اتوکلاو استفاده می‌شود؟'
$ws.Range("A320").Value = 'در سندرم دیسترس تنفسی حاد (ards) همه موارد زیر دیده می‌شوند به جز -'
$ws.Range("A321").Value = 'کدام یک از موارد زیر در نارسایی مزمن کلیه دیده نمی‌شود: سپتامبر 2010'
$ws.Range("A322").Value = 'مسدودکننده‌های h2 برای زخم دوازدهه به مدت چند وقت تجویز می‌شوند؟'
$ws.Range("A323").Value = 'سندرم میگ شامل همه موارد زیر به جز کدام است؟'
$ws.Range("A324").Value = 'یک کودک 4 ساله که برای کرانیوتومی آماده شده بود، لوله گذاری شد. ناگهان دستگاه بیهوشی شروع به نشان دادن باد کردن کیسه کرد. اقدام بعدی چیست؟'
$ws.Range("A325").Value = 'luxato erecta نام دیگری برای چیست؟'
$ws.Range("A326").Value = 'کودکی با حملات مکرر سرفه شدید همراه با صدای "ووپ" مراجعه می‌کند. کدام یک از گزینه‌های زیر بهترین نمونه برای جداسازی عامل بیماری و تأیید تشخیص محسوب می‌شود؟'
$ws.Range("A327").Value = 'یک پدر 45 ساله با علائم محرومیت از خواب، بیحالی، سردرد و خلق پایین دو ماه پس از اطلاع از ابتلای پسرش به لوسمی مراجعه میکند. او با دیگران به طور معقولانهای تعامل دارد، اما از کار خود غیبت کرده است. کدام یک از موارد زیر محتملترین تشخیص است؟'
$ws.Range("A328").Value = 'کدام مورد در پوست کراتینه شده استفاده می‌شود؟'
$ws.Range("A329").Value = 'همه موارد زیر از علل میوکاردیت هستند به جز'
$ws.Range("A330").Value = 'در مورد آسم ذاتی کدام گزینه صحیح است؟  
الف) شدیدتر و پایدارتر  
ب) ige طبیعی  
ج) سابقه خانوادگی مثبت  
د) تست پوستی مثبت'
$ws.Range("A331").Value = 'کدام یک از موارد زیر جزو قندهای کاهنده نیست؟'
$ws.Range("A332").Value = 'ظاهر "آسمان ستاره‌ای" در کدام مورد دیده می‌شود؟'
$ws.Range("A333").Value = 'در مورد سندرم کارتاژنر همه موارد زیر صحیح است به جز-'
$ws.Range("A334").Value = 'کاهش فشار اندولنفات در کدام یک از موارد زیر مشاهده می‌شود؟'
$ws.Range("A335").Value = 'یک آشپز ساندویچ‌هایی برای ۱۰ نفر که به پیک‌نیک می‌روند آماده می‌کند. هشت نفر از آن‌ها در عرض ۴ تا ۶ ساعت پس از مصرف ساندویچ‌ها دچار گاستروانتریت شدید می‌شوند. احتمالاً در بررسی‌ها مشخص می‌شود که آشپز ناقل کدام یک از موارد زیر است:'
$ws.Range("A336").Value = 'یک مرد 30 ساله با کلستئاتوم اتيك در گوش چپ و ترومبوفلبیت سینوس جانبی مراجعه کرده است. کدام یک از موارد زیر عمل انتخابی خواهد بود؟'
$ws.Range("A337").Value = 'جلوگیری از جابجایی رو به جلو استخوان درشت‌نی و ران توسط کدام ساختار انجام می‌شود؟'
$ws.Range("A338").Value = 'کلروما ناشی از کدام مورد است؟'
$ws.Range("A339").Value = 'فرزند پس از مرگ چیست؟'
$ws.Range("A340").Value = 'بزرگ‌ترین عصب جمجمه‌ای کدام است؟'
$ws.Range("A341").Value = 'آزمون رورشاخ با چه چیزی انجام می‌شود؟'
$ws.Range("A342").Value = 'طناب‌های شبکه بازویی بر اساس ارتباطشان با شریان زیربغلی که در پشت کدام عضله قرار دارد نامگذاری شده‌اند؟'
$ws.Range("A343").Value = 'درمان انتخابی در واکنش نوع یک جذام با نوریت شدید کدام است؟'
$ws.Range("A344").Value = 'کدام یک از ویروس‌های بدون پوشش و دارای rna تک رشته‌ای است؟'
$ws.Range("A345").Value = 'در مورد سندرم بندیکت، همه موارد زیر صحیح هستند به جز:'
$ws.Range("A346").Value = 'تمام موارد زیر در مورد پورپورای هنوخ شونلاین درست است به جز'
$ws.Range("A347").Value = 'همه موارد زیر در مورد تومورهای کارسینوئید پیش‌معده صحیح هستند، به جز:'
$ws.Range("A348").Value = 'تیتراسیون دوز یک دارو با پاسخ را می توان با کدام یک از راه های تجویز زیر انجام داد:'
$ws.Range("A349").Value = 'شریان قوس حلقی دوم چیست؟'
$ws.Range("A350").Value = 'یک مرد ۱۷ ساله با علائم جدید خستگی، بی‌حالی، تب و گلو درد به کلینیک مراجعه می‌کند. وی سابقه پزشکی قابل توجهی ندارد و تحت هیچ درمان دارویی نیست. معاینه فیزیکی به جز غدد لنفاوی گردنی بزرگ شده و قابل لمس، کاملاً طبیعی است. او کاهش وزن یا تعریق شبانه را گزارش نمی‌دهد. آزمایش‌های انجام شده شامل عکس قفسه سینه طبیعی، سوآب گلو منفی، اما فیلم خون غیرطبیعی با لنفوسیت‌های آتیپیک است. هموگلوبین ۱۵.۵ گرم در دسی‌لیتر؛ هماتوکریت ۴۲٪؛ پلاکت‌ها ۲۹۰۰۰۰ در میلی‌لیتر؛ گلبول سفید ۱۰۵۰۰ در میلی‌لیتر، با ۴۵٪ نوتروفیل‌های سگمانه، ۱٪ ائوزینوفیل‌ها و ۵۴٪ لنفوسیت‌ها که ۳۶٪ آن‌ها آتیپیک بودند. کدام یک از آزمایش‌های تشخیصی زیر مناسب‌ترین گزینه اولیه است؟'
$ws.Range("A351").Value = 'فشار خون به طور کاذب بالا در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A352").Value = 'عملکرد گوش میانی چیست؟'
$ws.Range("A353").Value = 'دوز آدرنالین تجویز شده برای یک کودک مبتلا به ایست قلبی چقدر است؟'
$ws.Range("A354").Value = 'عبارت صحیح در مورد هایپرترمی بدخیم کدام است؟'
$ws.Range("A355").Value = 'عفونت استافیلوکوکی باعث ایجاد تمام بیماری های پوستی زیر می شود به جز:-'
$ws.Range("A356").Value = 'یک کودک یک ساله با تحلیل عضلانی، از دست دادن چربی زیرپوستی بدون علائم ادم و وزن کمتر از ۶۰٪ استاندارد who مراجعه می‌کند. مادر سابقه عدم تغذیه کافی پروتئین و سایر مواد مغذی را پس از شش ماهگی کودک گزارش می‌دهد. تشخیص احتمالی چیست؟'
$ws.Range("A357").Value = 'مرد 59 ساله ای در تنفس از طریق بینی مشکل دارد. در معاینه، پزشک متوجه تورم غشاهای مخاطی در مئاتوس فوقانی بینی می شود. کدام مجرای سینوس های پارانازال به احتمال زیاد مسدود شده است؟'
$ws.Range("A358").Value = 'بیماری گرانولوماتوز مزمن به دلیل کمبود کدام مورد است؟'
$ws.Range("A359").Value = 'تحت شرایط فیزیولوژیکی طبیعی، فشار مایع مغزی-نخاعی (csf) با کدام یک از عوامل زیر متناسب است؟'
$ws.Range("A360").Value = 'تپش زبان کوچک با هر ضربان قلب در نارسایی شدید آئورت چه نامیده می‌شود؟'
$ws.Range("A361").Value = 'کدام گزینه در مورد هیستوپلاسما کپسولاتوم صحیح نیست؟'
$ws.Range("A362").Value = 'دیواره‌های پروگزیمال در آماده‌سازی دندان کلاس i برای آمالگام باید'
$ws.Range("A363").Value = 'عدم بینش (lack of insight) در کدام یک از موارد زیر دیده نمی‌شود؟'
$ws.Range("A364").Value = 'در تمامی موارد زیر به جز یکی، تلقیح مصنوعی با اسپرم همسر توصیه می‌شود.'
$ws.Range("A365").Value = 'علامت گاور در کدام یک از شرایط زیر کلاسیک است؟'
$ws.Range("A366").Value = 'کدام گزینه در مورد اریتم مولتیفرم صحیح نیست؟'
$ws.Range("A367").Value = 'عصب سه قلو در سیستم عصبی مرکزی چند هسته دارد؟'
$ws.Range("A368").Value = 'آلکالوز متابولیک در کدام مورد مشاهده می‌شود؟'
$ws.Range("A369").Value = 'آنستوموز پورتوكاوال در کدام محل مشاهده می‌شود؟'
$ws.Range("A370").Value = 'مخرج در نرخ بروز چیست؟'
$ws.Range("A371").Value = 'معمولاً تغییرات در یک رفتار مشکل‌زا چند مرحله را نشان می‌دهد:'
$ws.Range("A372").Value = 'کدام یک از مکانیسم های عمل سم خورنده نیست؟'
$ws.Range("A373").Value = 'کدام یک از عبارات زیر در مورد سینتیک مرتبه اول صحیح است؟'
$ws.Range("A374").Value = 'کدام نوروپاتی ارثی با چهارگانه کلاسیک نوروپاتی محیطی، رتینیت پیگمنتوزا، آتاکسی مخچه‌ای و پروتئین بالای مایع مغزی-نخاعی همراه است؟'
$ws.Range("A375").Value = 'کدام دندان شیری طولانی‌ترین مدت عملکرد را دارد؟'
$ws.Range("A376").Value = 'کدام ساختار مشخص شده در زیر مسئول گسترش عفونت به قسمت قدامی میانی بینی از طریق مجرای فرونتونازال است؟'
$ws.Range("A377").Value = 'تشخیص گره لنفانی سنتینل در حین عمل جراحی در ناحیه زیر بغل با استفاده از کدام روش انجام می‌شود؟'
$ws.Range("A378").Value = 'مرواریدهای عنبیه برای کدام بیماری پاتوگنومونیک هستند؟'
$ws.Range("A379").Value = 'اولین اسید چرب تشکیل شده در سنتز اسید چرب کدام است؟'
$ws.Range("A380").Value = 'آرایش سیم مرغی (chicken wire) مربوط به کدام نوع فیبر کلاژن است؟'
$ws.Range("A381").Value = 'کدام یک از انواع توهم‌های زیر برای اسکیزوفرنی پاتوگنومونیک است؟'
$ws.Range("A382").Value = 'جسم خارجی مسی در چشم، کدام یک از موارد زیر مشاهده خواهد شد؟'
$ws.Range("A383").Value = 'میوپاتی ناشی از دارو می‌تواند توسط تمام موارد زیر ایجاد شود به جز:'
$ws.Range("A384").Value = 'کدام یک از موارد زیر از ویژگی‌های سیستم ایمنی مخاطی است؟'
$ws.Range("A385").Value = 'میکسوویروس‌ها شامل کدام‌یک از موارد زیر هستند؟  
الف) اورتومیکسوویروس  
ب) آنفلوانزا  
ج) سرخک  
د) فلج اطفال  
ه) ویروس هرپس سیمپلکس'
$ws.Range("A386").Value = 'آنژین روده ای مجموعه علائم زیر است-'
$ws.Range("A387").Value = 'سدیم فلوراید با مهار کدام آنزیم، گلیکولیز را مهار می‌کند؟'
$ws.Range("A388").Value = 'پرونئوس برویس به کدام قسمت متصل می‌شود؟'
$ws.Range("A389").Value = 'آسپرین چند روز قبل از عمل جراحی قطع می‌شود؟'
$ws.Range("A390").Value = 'در مورد سندرم داندی-واکر، همه موارد زیر دیده می‌شوند به جز-'
$ws.Range("A391").Value = 'پس از حمله انفارکتوس میوکارد، مرگ و میر و عوارض بیمار با کدام مورد مشخص می‌شود؟'
$ws.Range("A392").Value = 'کدام یک از موارد زیر تومور نوروژنیک اولیه است؟'
$ws.Range("A393").Value = 'کدام آنتی‌بیوتیک بیشترین نقش را در ایجاد dili (آسیب کبدی ناشی از دارو) دارد؟'
$ws.Range("A394").Value = 'کدام یک از موارد زیر طولانی‌ترین اثر را در میان مسدودکننده‌های بتای چشمی دارد؟'
$ws.Range("A395").Value = 'درمان کزاز کدام است؟'
$ws.Range("A396").Value = 'همه موارد زیر در مورد پروتئین های دناتوره شده صحیح است به جز:'
$ws.Range("A397").Value = 'نوع بالینی اسکلریت مرتبط با بیماری‌های کلاژن کدام است؟'
$ws.Range("A398").Value = 'نوزاد متولد شده از مادر دیابتی همه موارد زیر را دارد به جز:'
$ws.Range("A399").Value = 'حس چشایی از یک سوم خلفی زبان توسط کدام عصب منتقل می‌شود؟'
$ws.Range("A400").Value = 'مخزن باونلا آنژیوماتوزیس چیست؟'
$ws.Range("A401").Value = 'پادزهر مسمومیت با استریکنین کدام است؟'
$ws.Range("A402").Value = 'یک کودک 10 ساله با سینوس پیش گوشی مراجعه می‌کند. هیچ ترشح یا التهابی وجود ندارد. سینوس مشابهی نیز در مادرش مشاهده شده بود که در سن 30 سالگی تحت درمان قرار گرفت. خط مدیریت شما چه خواهد بود؟'
$ws.Range("A403").Value = 'یک کودک 3 ساله با عفونت دستگاه تنفسی فوقانی به بخش سرپایی آورده می‌شود. با پرسش درباره سابقه واکسیناسیون قبلی، مشخص می‌شود که کودک از بدو تولد هیچ واکسیناسیون اولیه‌ای دریافت نکرده است. مناسب‌ترین اقدام در مدیریت این بیمار چیست؟'
$ws.Range("A404").Value = 'شایع ترین عارضه بعدی اسپوندیلیت سلی در یک نوجوان کدام است؟'
$ws.Range("A405").Value = 'تفاوت بین هیوسین و آتروپین این است که هیوسین:'
$ws.Range("A406").Value = 'تتراد اورتریت، آرتریت، کنژنکتیویت و زخم‌های دهانی در کدام مورد مشاهده می‌شود؟'
$ws.Range("A407").Value = 'یک زن ۲۵ ساله که ۲ سال پیش آییودی گذاشته است، با تب مراجعه می‌کند. در معاینه، actinomyces در سیتولوژی دهانه رحم مثبت است. در مورد آییودی چه توصیه‌ای می‌کنید؟'
$ws.Range("A408").Value = 'شکل ذخیره‌شده هورمون تیروئید -'
$ws.Range("A409").Value = 'کاهش حجم ریه یکی از ویژگی‌های کدام بیماری است؟'
$ws.Range("A410").Value = 'مهم‌ترین عامل پیش‌آگهی در تومور ویلمز-'
$ws.Range("A411").Value = 'پیودرمای گانگرنوزوم در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A412").Value = 'شاخه‌های ورید کلیوی چپ همه موارد زیر هستند به جز -'
$ws.Range("A413").Value = 'آنتی‌ژن‌های پلی‌ساکاریدی:'
$ws.Range("A414").Value = 'یک ویژگی خاص فیبروبلاست های رباط پریودنتال چیست؟'
$ws.Range("A415").Value = 'آوران رفلکس کرماستریک شامل کدام مورد است؟'
$ws.Range("A416").Value = 'نوار قلب (ecg) در تشخیص ایسکمی در مناطقی که توسط کدام یک از عروق زیر خونرسانی می‌شوند، ضعیف عمل می‌کند؟'
$ws.Range("A417").Value = 'پایش درمانی غلظت پلاسمایی داروهای ضد فشار خون انجام نمی‌شود زیرا:'
$ws.Range("A418").Value = 'گره هبردن نشان دهنده درگیری کدام مفصل است؟'
$ws.Range("A419").Value = 'آلوپسی موضعی (شبیه به خوردگی توسط بید) در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A420").Value = 'در کدام یک از موارد زیر در مورد سرطان پستان، پرتودرمانی پس از ماستکتومی مورد نیاز نیست؟'
$ws.Range("A421").Value = 'هیستوگرام به عنوان روشی برای نمایش گروهی داده‌های زیر استفاده می‌شود:'
$ws.Range("A422").Value = 'نیمه عمر t3'
$ws.Range("A423").Value = 'خود (ایگو) توسط چه چیزی اداره می‌شود؟'
$ws.Range("A424").Value = 'تمام موارد زیر در مورد مواد افیونی صحیح است به جز؟'
$ws.Range("A425").Value = 'زنی که در هنگام هرس کردن بوته‌های گل رز انگشت خود را سوراخ کرده است، دچار یک پوسچول موضعی می‌شود که به زخم تبدیل می‌شود. سپس چندین ندول در مسیر تخلیه لنفاوی موضعی ایجاد می‌شود. محتمل‌ترین عامل اتیولوژیک کدام است؟'
$ws.Range("A426").Value = 'ضخامت باند ماتریکس کلاس ii چقدر است؟'
$ws.Range("A427").Value = 'اسکلرای آبی فیزیولوژیک در چه مواردی دیده می‌شود؟'
$ws.Range("A428").Value = 'گال در بزرگسالان شامل کدام ناحیه نمی‌شود؟'
$ws.Range("A429").Value = 'داروی انتخابی برای طاعون کدام است؟'
$ws.Range("A430").Value = 'اسکوتولا به طور کلاسیک در کدام مورد دیده می‌شود؟'
$ws.Range("A431").Value = 'کاردیومگالی کشان -'
$ws.Range("A432").Value = 'شایع ترین ناهنجاری سیتوژنتیک در سندرم میلودیسپلاستیک (mds) بالغین کدام است؟'
$ws.Range("A433").Value = 'در مقایسه با شیر گاو، شیر انسان بیشتر حاوی کدام ماده است؟'
$ws.Range("A434").Value = 'کدام یک از عبارات زیر در مورد لپتوسپیروز صحیح است؟'
$ws.Range("A435").Value = 'یک دختر ۵ ساله با پورپورای قابل لمس روی باسن، آرترالژی، درد شکم همراه با اسهال و دفع خون از رکتوم مراجعه کرده است. همچنین بیمار پروتئینوری دارد. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A436").Value = 'یک دختر 3 ساله با عفونت ادراری مکرر مراجعه می‌کند. در سونوگرافی هیدرونفروز با نقص پر شدن و سایه منفی مثانه بدون دهانه نابجا مشاهده می‌شود.'
$ws.Range("A437").Value = 'کدام یک از رنگ‌آمیزی‌های ویژه زیر می‌تواند برای تشخیص رابدومیوسارکوم استفاده شود؟'
$ws.Range("A438").Value = 'دنده کاذب کدام است؟'
$ws.Range("A439").Value = 'کدام یک از موارد زیر یک اختلال اتوزومی غالب است؟'
$ws.Range("A440").Value = 'همه موارد زیر از عوامل القای بیهوشی داخل وریدی هستند به جز:'
$ws.Range("A441").Value = 'چه زمانی می‌توان دیسترس تنفسی حاد را در یک نوزاد تشخیص داد؟'
$ws.Range("A442").Value = 'پروبنسید با کدام یک از موارد زیر تداخل دارد؟'
$ws.Range("A443").Value = 'سندرم کلاینفلتر با کدام روش تشخیص داده می‌شود؟'
$ws.Range("A444").Value = 'ظاهر مدرسه ماهی توسط کدام یک نشان داده می‌شود؟'
$ws.Range("A445").Value = 'همه موارد زیر در مورد عصب فک پایین صحیح هستند به جز:'
$ws.Range("A446").Value = 'کدام یک از موارد زیر مشتق از لوله‌های مزونفریک در مردان است؟'
$ws.Range("A447").Value = 'کدام یک از داروهای ضد مالاریا در دوران بارداری ایمن‌تر است؟'
$ws.Range("A448").Value = 'نظر خود را در مورد تشخیص نمودار ecg ارائه شده بیان کنید.'
$ws.Range("A449").Value = 'احتمال ابتلای یک کارمند بهداشتی به hiv در اثر فرو رفتن تصادفی سوزن چقدر است؟'
$ws.Range("A450").Value = 'عارضه شایع آبله مرغان در کودکان کدام است؟'
$ws.Range("A451").Value = 'بیمار مبتلا به سیروز، hbs ag+ مثبت و سطح افزایش یافته آلفا فتوپروتئین دارد. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A452").Value = 'کدام یک از موارد زیر ویژگی یا خاصیت یک ویروس نیست؟'
$ws.Range("A453").Value = 'بافت‌ای که بیشترین حساسیت را به پرتو دارد کدام است؟'
$ws.Range("A454").Value = 'گاز خون شریانی یک بیمار کاهش ph، افزایش pco2 و بیکربنات بالا را نشان می‌دهد. تشخیص:'
$ws.Range("A455").Value = 'چه زمانی نسبت مایع داخل سلولی (icf) و مایع خارج سلولی (ecf) در کودک به سطح بزرگسال می‌رسد؟'
$ws.Range("A456").Value = 'تب بخار فلزی با کدام یک از موارد زیر مرتبط است؟'
$ws.Range("A457").Value = 'کمبود کدام آنزیم باعث عدم تحمل ارثی فروکتوز می‌شود؟'
$ws.Range("A458").Value = 'همه موارد زیر از ویژگی‌های نوزاد نارس هستند، به جز –'
$ws.Range("A459").Value = 'کدام یک از موارد زیر معمولاً اولین نشانه بلوغ در دختران است؟'
$ws.Range("A460").Value = 'آبسه لگنی با کدام یک از موارد زیر تظاهر می‌کند؟'
$ws.Range("A461").Value = 'فردی در اثر اصابت توپ تنیس به چشم راست دچار آسیب شده است. پس از آن درد ایجاد شده و در معاینه ته چشم، یک نقطه قرمز در ماکولا دیده می‌شود. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A462").Value = 'هماتوم اپیدورال با چه درصدی از تروماهای شدید همراه است؟'
$ws.Range("A463").Value = 'مایعات ترانس سلولی در کجا وجود دارند؟'
$ws.Range("A464").Value = 'افزایش سطح هموگلوبین a2 مشخصه کدام یک از موارد زیر است؟'
$ws.Range("A465").Value = 'یک زن ۲۸ ساله با سابقه چندساله اسهال متناوب و درد شکمی برای بیماری التهابی روده بررسی می‌شود. ارزیابی آندوسکوپی ایلئوم انتهایی، کولون و رکتوم انجام شده است. کدام یک از مشاهدات آندوسکوپی زیر بیشتر نشان‌دهنده بیماری کرون است تا کولیت اولسراتیو؟'
$ws.Range("A466").Value = 'مرز چپ قلب در رادیوگرافی قفسه سینه (cxr) توسط کدام ساختار تشکیل شده است؟ الف) شریان ریوی ب) ورید ریوی ج) آئورت شکمی د) قوس آئورت ه) بطن راست'
$ws.Range("A467").Value = 'یک مرد 23 ساله به دلیل علائم اختلالات خواب و ''احساس افسردگی'' به درمانگاه روانپزشکی مراجعه می‌کند. علائم دو ماه پیش شروع شد، زمانی که متوجه شد دوست‌دخترش به او خیانت کرده و با رئیسش در رابطه بوده است. بیمار بلافاصله رابطه را قطع کرد. اما از آن زمان احساس غمگینی دارد و گزارش می‌دهد که هیچ چیز دیگر او را خوشحال نمی‌کند. بیمار دیگر به محل کار نمی‌رود و حتی از اتاقش خارج نمی‌شود مگر در موارد کاملاً ضروری. او به مدت هفت روز نخوابیده و در دو ماه گذشته چند کیلوگرم وزن از دست داده است. در طول مصاحبه گفت: ''من یک بازنده هستم، دوست‌دخترم به جای من یک مرد 50 ساله را انتخاب کرد، حتماً مشکلی جدی در من وجود دارد.'' تشخیص چیست؟'
$ws.Range("A468").Value = 'ارثی اسیدوری اوروتیک نوع اول به دلیل کمبود کدام آنزیم است؟'
$ws.Range("A469").Value = 'بیشترین پتانسیل سایش در کدام حالت وجود دارد؟'
$ws.Range("A470").Value = 'عمل سوئنسون برای کدام مورد انجام می‌شود؟'
$ws.Range("A471").Value = 'طحال بسیار بزرگ شده همراه با پان سیتوپنی -'
$ws.Range("A472").Value = 'کدام گزینه در مورد آدنومیوز صحیح است؟'
$ws.Range("A473").Value = 'در حمایت پایه حیات (bls)، حمایت به کدام یک از اعضای زیر ارائه می‌شود؟'
$ws.Range("A474").Value = 'یک پسر ۱۰ ساله با احساس ناراحتی در پشت جناغ به بیمارستان مراجعه می‌کند. سی‌تی اسکن یک تومور میانی در غده تیموس را نشان می‌دهد. کدام یک از رگ‌های زیر به احتمال زیاد توسط این تومور فشرده می‌شود؟'
$ws.Range("A475").Value = 'آسپرین برای پیشگیری ثانویه از بیماری ایسکمیک قلب استفاده می‌شود زیرا:'
$ws.Range("A476").Value = 'تمام موارد زیر در مورد سندرم کاگنر صحیح است، به جز'
$ws.Range("A477").Value = 'موارد منع استفاده از ترمیم‌های ریخته‌گری (کست) همه موارد زیر هستند به جز:'
$ws.Range("A478").Value = 'انقباض و انبساط متناوب مردمک‌ها که به عنوان علامت هیپوس شناخته می‌شود، در کدام مورد مشاهده می‌شود؟'
$ws.Range("A479").Value = 'چه چیزی از دیافراگم در سطح t12 عبور می‌کند؟'
$ws.Range("A480").Value = 'یک زن ۲۳ ساله تحت تیروئیدکتومی کامل به دلیل کارسینوم غده تیروئید قرار می‌گیرد. در روز دوم پس از عمل، او از احساس سوزن سوزن شدن در دستانش شکایت می‌کند. او به نظر بسیار مضطرب می‌رسد و بعداً از گرفتگی عضلات شکایت می‌کند. کدام یک از موارد زیر مناسب‌ترین استراتژی مدیریت اولیه است؟'
$ws.Range("A481").Value = 'کدام یک از موارد زیر منشأ ویروسی ندارد؟'
$ws.Range("A482").Value = 'همه موارد زیر در سندرم حذف 22q11.2 دیده می‌شوند به جز'
$ws.Range("A483").Value = 'بهترین شاخص برای تعیین حداکثر منفعت برای جامعه از طریق راهبردهای مداخله پیشگیرانه کدام است؟'
$ws.Range("A484").Value = 'در یک فرد بالغ 50 کیلوگرمی، چه مقدار مایع باید در 8 ساعت اول برای سوختگی درجه دوم 40٪ تجویز شود؟'
$ws.Range("A485").Value = 'فردی که دارای صفت هتروزیگوت سلول داسی شکل است در برابر عفونت کدام یک محافظت می شود؟'
$ws.Range("A486").Value = 'یک پسر 10 ساله پشتش را سوزانده است. این سوختگی چه درصدی از سطح بدن او را شامل می‌شود؟'
$ws.Range("A487").Value = 'تشریح آنزیم های غیرفعال کننده، مکانیسم مهم مقاومت دارویی در بین تمام این آنتی بیوتیک ها به جز کدام یک است؟'
$ws.Range("A488").Value = 'واکسن هاری اولین بار توسط چه کسی توسعه یافت؟'
$ws.Range("A489").Value = 'کدام یک از موارد زیر از حلقه ناف اولیه عبور نمی‌کند؟'
$ws.Range("A490").Value = 'داروی انتخابی برای درمان عفونت ناشی از استافیلوکوک اورئوس مقاوم به متی سیلین کدام است؟'
$ws.Range("A491").Value = 'شایع‌ترین علت تومور کروکنبرگ چیست؟'
$ws.Range("A492").Value = 'غلات فاقد کدام ماده هستند؟'
$ws.Range("A493").Value = 'یک بیمار مبتلا به هپاتیت مزمن با نتایج مثبت hbsag و anti-hbc (igg) همچنین تشخیص عفونت hiv را دریافت کرده است. بهترین دارو برای حمله به هر دو ویروس کدام است؟'
$ws.Range("A494").Value = 'کم خونی همولیتیک ممکن است با همه موارد زیر به جز ... مشخص شود.'
$ws.Range("A495").Value = 'کدام یک از موارد زیر با سندرم داون مرتبط است: مارس 2010'
$ws.Range("A496").Value = 'رفلکس چشم گربه‌ای آموروتیک در کدام مورد مشاهده می‌شود؟'
$ws.Range("A497").Value = 'رایج‌ترین شکل تنوع dna چیست؟'
$ws.Range("A498").Value = 'کدام یک از موارد زیر در مورد کاهش ضربان قلب زودرس در کاردیوتوکوگرافی (ctg) صحیح است؟'
$ws.Range("A499").Value = 'خطرناک‌ترین عارضه در بیمارانی که تحت عمل جراحی تیروئید قرار گرفته‌اند و در محل عمل هماتوم ایجاد شده است چیست؟'
